# Generate Report for Handoff
# Updates the localization-status workbook to reflect a fresh handoff
# generation run: the "Latest HO Xliff Generate Date" / "Latest Handoff
# Datetime" timestamps for the four "Ready for handoff" files move
# forward, and their Priority is normalized from "low" to "ht".

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

for ($r = 4; $r -le 7; $r++) {
    # Overview sheet: "Latest HO Xliff Generate Date" column (G)
    $wsOverview.Range("G$r").Value = "2016-09-06 18:42:37"

    # zh-cn sheet: Priority (E) and Latest Handoff Datetime (H)
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-09-06 18:42:32"

    # de-de sheet: Priority (E) and Latest Handoff Datetime (H)
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-09-06 18:42:37"
}
